$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.102.70'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").Value = '1.832.45'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9986'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6619'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2945'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07322'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.86%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07644'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("D12").Value = '1.838.76'
$ws.Range("E12").Value = '  -0.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.008'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6727'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '85.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.122'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.71%  '
$ws.Range("D17").Value = '29.082.73'
$ws.Range("E17").Value = '  -1.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008189'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '227.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.48%  '
$ws.Range("E20").Value = '  -2.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.244'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9992'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1418'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.639'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.499'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.219'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.102'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.198'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05307'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7462'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.847'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.127'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.679'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("D37").Value = '1.297.01'
$ws.Range("E37").Value = '  -2.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01809'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.702'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("E40").Value = '  -4.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.013'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9983'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.94%  '
$ws.Range("D44").Value = '1.984.56'
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5174'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000121'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.747'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '63.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.206'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07368'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05908'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.31%  '
